# Kanban.xlsx update - v.1.4
# Reworks the "base de datos" (database) tasks:
#  - the old To-do item "Diseña el esquema inicial de la base de datos
#    (SQLite/PostgreSQL)" is removed from the To do column,
#  - a shortened version of it ("Diseña el esquema inicial de la base de
#    datos") is inserted into the Done column, right after
#    "Configurar las URLs para la API",
#  - three brand-new Done items are appended at the end of the list,
#  - the remaining To do item shifts down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1) Grow the grid: copy the blank-row formatting (row 23, which holds the
#    "empty A/B/C" look used for the bottom rows) down across the 3 new
#    rows (24-26) that the sheet needs.
$ws.Range("A23:C23").Copy()
$ws.Range("A24:C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Give the new Done entry that lands on row 20 the same "To do" style
#    (s=1) that the old column-A task used, and give the other new/shifted
#    Done rows (21-24) the regular continuation style (s=5) that the rest
#    of the Done column uses.
$ws.Range("A21").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C11").Copy()
$ws.Range("C21:C24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Write the final text for every affected cell. The three brand-new
#    Done items are entered first (in the order they were authored), then
#    the shortened "esquema inicial" task is written last, matching the
#    shared-string order of the saved workbook.
$ws.Range("C21").Value = "Probar los endpoints en Postman o el navegador"
$ws.Range("C22").Value = "Definir el esquema de la base de datos"
$ws.Range("C23").Value = "Crear las migraciones y aplicar cambios en la base de datos"
$ws.Range("C24").Value = "Verificar el esquema en Django Admin y probarlo"
$ws.Range("C20").Value = "Diseña el esquema inicial de la base de datos"

# The remaining To do task moves from A22 down to A25; clear the old spot
# and write it into its new place.
$ws.Range("A21").ClearContents()
$ws.Range("A22").ClearContents()
$ws.Range("A25").Value = "Diseña una página básica con React y conecta al backend"

# 4) Resize the Excel table ("Tabla1") so it covers the new extent and
#    refresh the autofilter/dimension along with it.
$tbl.Resize($ws.Range("A1:C26"))

# 5) Match the view state captured in the saved workbook.
$ws.Range("A22").Select()
$excel.ActiveWindow.ScrollRow = 6
